# Fruta / hortaliza, semanal
#
# A new weekly price record was inserted into the "Apio" sheet at row 132
# (Feria Lagunitas de Puerto Montt, "Primera" quality, date 2021-11-05),
# pushing the previous rows 132-155 down to 133-156.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 132; this shifts the existing rows 132:155
# down to 133:156 (and grows the used range to A1:R156).
$ws.Rows.Item(132).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A132").Value = 4
$ws.Range("B132").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C132").Value = "Los Lagos"
$ws.Range("D132").Value = 44505
$ws.Range("E132").Value = 10
$ws.Range("F132").Value = 100112017
$ws.Range("G132").Value = "Apio"
$ws.Range("H132").Value = "Americana (o)"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 50
$ws.Range("K132").Value = 9000
$ws.Range("L132").Value = 10000
$ws.Range("M132").Value = 9500
$ws.Range("N132").Value = "`$/docena de matas"
$ws.Range("O132").Value = "Región de Coquimbo"
$ws.Range("P132").Value = 1583
$ws.Range("Q132").Value = 6
$ws.Range("R132").Value = "Hortaliza"
